$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (domingo)
$ws.Range("C2").Value = ""
$ws.Range("E2").Value = ""

# Row 3 (segunda)
$ws.Range("C3").Value = ""
$ws.Range("E3").Value = ""

# Row 4 (terça)
$ws.Range("B4").Value = "ratatouille"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "linguiça"
$ws.Range("E4").Value = "abobrinha"

# Row 5 (quarta)
$ws.Range("B5").Value = "ratatouille"
$ws.Range("C5").Value = ""
$ws.Range("E5").Value = ""

# Row 6 (quinta)
$ws.Range("B6").Value = "ratatouille"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "bife de frango"
$ws.Range("E6").Value = "cenoura ralada"

# Row 7 (sexta)
$ws.Range("B7").Value = "ratatouille"
$ws.Range("C7").Value = ""
$ws.Range("E7").Value = ""

# Row 8 (sábado)
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = "carne moída"
$ws.Range("E8").Value = "maionese"
